$d = $word.ActiveDocument

# Helper: replace an entire paragraph's content with a single clean run of
# text (no stray w:proofErr / split runs). Works by inserting a brand-new
# paragraph (with the desired text) immediately before the target paragraph,
# then deleting the now-shifted original paragraph (including its paragraph
# mark) in one go -- a whole-paragraph delete removes any zero-width marks
# (like w:proofErr) that live inside it, unlike a partial in-paragraph delete.
function Clean-Paragraph($doc, $index, $newText) {
    $p = $doc.Paragraphs.Item($index)
    $full = $p.Range
    $full.InsertBefore($newText + "`r")
    $dirty = $doc.Paragraphs.Item($index + 1)
    $dirty.Range.Delete()
}

# 1) "This is my World Quant University ... Vietnam stock market" paragraph:
#    merge the two split runs ("...stock " + "market") back into one run.
Clean-Paragraph $d 1 "This is my World Quant University Capstone Project: Regime Change Detection by applying a Directional-change Event approach in the Vietnam stock market"

# 2) "The idea is to use to Directional Change Event ... FiinPro." paragraph:
#    merge the split runs ("...collected from " + "FiinPro" + ".") into one run.
Clean-Paragraph $d 5 "The idea is to use to Directional Change Event approach to detect Regime Change in market movement. The target is the Vietnam Stock Index - VNINDEX, which the historical data is collected from FiinPro."

# 3) "The code to calculate the DC Log Return Indicator ..." paragraph:
#    merge the split runs ("...Thomas " + "github" + " code: " + url) into one run.
Clean-Paragraph $d 7 "The code to calculate the DC Log Return Indicator is modified from Thomas github code: https://github.com/ThomasWangWeiHong/Time-Series-Directional-Change-Analysis"

# 4) "We them fit the DC Log Return ... found here: <hyperlink>" paragraph:
#    merge every split text run before the hyperlink into one run, while
#    leaving the hyperlink field/run completely untouched.
$p9 = $d.Paragraphs.Item(9)
$hyperlink = $d.Hyperlinks.Item(1)

# Split the paragraph into two paragraphs right before the hyperlink so the
# text-only portion becomes its own whole paragraph.
$splitPoint = $d.Range($hyperlink.Range.Start, $hyperlink.Range.Start)
$splitPoint.InsertBefore("`r")

# Clean that text-only paragraph the same way as above.
Clean-Paragraph $d 9 "We them fit the DC Log Return and Normal Log Return into Hidden Markov Model(HMM) to detect hidden state, then we can identify the normal and abnormal regimes. The formula of can be found here: "

# Re-merge the cleaned text paragraph with the following hyperlink paragraph
# by deleting the paragraph mark that now separates them.
$p9b = $d.Paragraphs.Item(9)
$markRange = $d.Range($p9b.Range.End - 1, $p9b.Range.End)
$markRange.Delete()

# 5) Add a new paragraph right after that paragraph (the one ending with the
#    hyperlink), before the following blank paragraph:
#    "We also try impact of different thetas, you can find the result in the
#    different theta folders."
$p9c = $d.Paragraphs.Item(9)
$endOfP9 = $d.Range($p9c.Range.End - 1, $p9c.Range.End - 1)
$endOfP9.InsertAfter("`rWe also try impact of different thetas, you can find the result in the different theta folders.")
